# Backlog.xlsx edit:
#  - Column C ("Semana") on both sheets changes from the text "Semana 01"
#    to the numeric value 1 for every data row.
#  - The active sheet switches from SPN to ITI, and the selection / scroll
#    position on each sheet is updated to match the new state.

$wb = $excel.ActiveWorkbook

# --- SPN sheet (first sheet) ---------------------------------------------
$spn = $wb.Worksheets.Item("SPN")
$spn.Range("C2:C26").Value = 1

# --- ITI sheet (second sheet) ---------------------------------------------
# Column C picks up the plain "centered, no fill" style (same as SPN's column
# C) now that it holds a number instead of a shared string.
$iti = $wb.Worksheets.Item("ITI")
$iti.Range("C2:C10").Value = 1
$iti.Range("C2:C10").Interior.Pattern = -4142
$iti.Range("C2:C10").HorizontalAlignment = -4108

# --- View state: SPN is no longer the active tab, scrolled back to A1,
#     selection moved to C2 ------------------------------------------------
$spn.Activate()
$spn.Range("A1").Select()
$spn.Range("C2").Select()

# --- View state: ITI becomes the active tab, selection at D24 -------------
$iti.Activate()
$iti.Range("D24").Select()
